# Add a "Server 3.4.0" results worksheet, cloned from the existing
# "Server 3.3.0" sheet (same labels/formulas/styles) but with the measured
# run data cleared out, ready for the next round of numbers to be filled in.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate "Server 3.3.0", inserting the copy immediately after it.
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)

# The newly-inserted copy becomes the active sheet.
$ws2 = $wb.ActiveSheet
$ws2.Name = "Server 3.4.0"

# Clear the per-run measurement cells (the Min/Max/Average formulas and all
# row/column labels stay in place, so they recompute to 0 / #DIV/0!).
$ws2.Range("F2:L4").ClearContents()
$ws2.Range("F7:I9").ClearContents()

# Leave the selection on the first data cell of the fresh sheet.
$ws2.Range("F2").Select() | Out-Null
